$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (pushes existing rows 3-8 down to 4-9).
# The new row gets a copy of the (former) row 3's data for this market/product,
# but with an updated Fecha (date) and Volumen value, reflecting a new
# weekly price record for "Feria Lagunitas de Puerto Montt - Arándano (blue)".
$ws.Rows("3:3").Insert()

$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C3").Value = "Los Lagos"
$ws.Range("D3").Value = 44519
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100101
$ws.Range("H3").Value = "Berries"
$ws.Range("I3").Value = 100101001
$ws.Range("J3").Value = "Arándano (blue)"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 3700
$ws.Range("O3").Value = 3800
$ws.Range("P3").Value = 3750
$ws.Range("Q3").Value = "`$/kilo"
$ws.Range("R3").Value = "Región del Maule"
$ws.Range("S3").Value = 3750
$ws.Range("T3").Value = 1
